# "Comunicacao Cross e Painel"
#
# - Translate the patient names in column B (Name) from English
#   ("Patient N") to Portuguese ("Paciente N").
# - Update the PLA_Id for reservations 3 and 5 (rows 4 and 6) to the
#   other plan id already used elsewhere in the sheet.
# - Correct the start/end time of reservation 6 (row 7) by pushing it
#   one hour later.
# - Re-fit the columns whose widest content changed (Name, PAC_Id,
#   PLA_Id) and restore the active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B ("Name"): Patient N -> Paciente N -------------------------
$ws.Range("B2").Value = "Paciente 1"
$ws.Range("B3").Value = "Paciente 2"
$ws.Range("B4").Value = "Paciente 3"
$ws.Range("B5").Value = "Paciente 4"
$ws.Range("B6").Value = "Paciente 5"
$ws.Range("B7").Value = "Paciente 6"
$ws.Range("B8").Value = "Paciente 7"

# --- Column G ("PLA_Id") value correction on rows 4 and 6 ---------------
$ws.Range("G4").Value = "5aef92b8124bc3b4db6932b7"
$ws.Range("G6").Value = "5aef92b8124bc3b4db6932b7"

# --- Row 7: reservation start/end time moved one hour later -------------
$ws.Range("C7").Value = 43266.958333333336
$ws.Range("D7").Value = 43267.041666666664

# --- Re-fit columns whose widest content changed -------------------------
$ws.Columns("B").EntireColumn.AutoFit()
$ws.Columns("E").EntireColumn.AutoFit()
$ws.Columns("G").EntireColumn.AutoFit()

# --- Restore the active-cell selection -----------------------------------
[void]$ws.Range("E12").Select()
